$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $valA = $cellA.Value()
    if ($valA -eq "🟥") { $cellA.Value = "📕" }
    elseif ($valA -eq "⬛") { $cellA.Value = "📘" }
    elseif ($valA -eq "🟧") { $cellA.Value = "📙" }
    elseif ($valA -eq "🟩") { $cellA.Value = "📗" }

    $cellB = $ws.Cells.Item($r, 2)
    $valB = $cellB.Value()
    if ($valB -eq "noir") { $cellB.Value = "bleu" }
}
